$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '61.040.62'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + '  +0.14%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'" + '3.392.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + '  -0.56%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'" + '571.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + '  -0.03%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'" + '142.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + '  +0.52%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'" + '  +0.00%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'" + '  -0.06%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'" + '7.64'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'" + '  +0.47%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'" + '0.124'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + '  -1.19%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'" + '0.391'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + '  -0.09%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'" + '3.971.19'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'" + '  -0.57%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'" + '  +1.82%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'" + '27.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + '  -1.33%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'" + '  +0.03%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'" + '3.395.04'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + '  -0.60%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'" + '61.147.32'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'" + '  +0.10%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'" + '  -3.09%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'" + '  -5.01%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'" + '8.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + '  -4.19%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'" + '382.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'" + '  -1.51%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'" + '74.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + '  +3.02%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'" + '0.554'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + '  -2.52%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'" + '  +0.43%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'" + '  -5.15%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'" + '3.523.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'" + '  -0.83%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'" + '0.181'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'" + '  +0.32%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'" + '  +0.07%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'" + '7.35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'" + '  -1.20%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'" + '8.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'" + '  -1.65%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'" + '  -0.58%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'" + '  -3.38%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'" + '  +0.02%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'" + '23.35'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'" + '  -2.23%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'" + '  -0.44%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'" + '166.58'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'" + '  -0.69%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'" + '  -1.63%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'" + '3.422.87'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'" + '  -0.49%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'" + '  -4.14%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'" + '  -1.74%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'" + '27.24'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'" + '  +0.87%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'" + 'Mantle'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'" + 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'" + '0.780'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + '  -1.67%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'" + 'FirstDigitalUSD'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'" + 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'" + '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'" + '  -0.07%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'" + '4.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'" + '  -2.14%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'" + '1.68'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + '  -1.83%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'" + '1.14'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'" + '  +0.12%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'" + '2.455.89'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'" + '  -5.18%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'" + '23.08'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'" + '  +0.72%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'" + '6.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'" + '  -3.30%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'" + 'dogwifhat'
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'" + 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'" + '2.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'" + '  +7.36%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'" + 'VeChain'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'" + 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'" + '0.0266'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + '  +1.93%  '
$ws.Range("E51").Style = "Normal"
